$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2.860461381038202
$ws.Range("D2").Value = 15.99302765733131
$ws.Range("E2").Value = 2.984030474818697
$ws.Range("F2").Value = 2.050004649001905
$ws.Range("G2").Value = 3.682990585644517
$ws.Range("H2").Value = 27.57051474783463
$ws.Range("C3").Value = 4.514746097148066
$ws.Range("D3").Value = 16.53916365872576
$ws.Range("E3").Value = 1.913302377929095
$ws.Range("F3").Value = 4.413063653873156
$ws.Range("G3").Value = 4.754242873000782
$ws.Range("H3").Value = 32.13451866067686
$ws.Range("C4").Value = 7.548704016121698
$ws.Range("D4").Value = 16.07401222490661
$ws.Range("E4").Value = 1.272007460608741
$ws.Range("F4").Value = 6.275224935914857
$ws.Range("G4").Value = 4.704909838082013
$ws.Range("H4").Value = 35.87485847563391
$ws.Range("C5").Value = 2.108977561824382
$ws.Range("D5").Value = 15.89871350328061
$ws.Range("E5").Value = 3.131533570470598
$ws.Range("F5").Value = 1.912260144828285
$ws.Range("G5").Value = 3.727828575685083
$ws.Range("H5").Value = 26.77931335608895
$ws.Range("C6").Value = 3.373517242612782
$ws.Range("D6").Value = 16.97781132775818
$ws.Range("E6").Value = 2.259896794268828
$ws.Range("F6").Value = 3.684442760960605
$ws.Range("G6").Value = 4.651326836305295
$ws.Range("H6").Value = 30.94699496190569
$ws.Range("C7").Value = 5.76163892911005
$ws.Range("D7").Value = 16.8267756289517
$ws.Range("E7").Value = 1.435972371366783
$ws.Range("F7").Value = 5.554879653671966
$ws.Range("G7").Value = 5.136801510576634
$ws.Range("H7").Value = 34.71606809367714
$ws.Range("C8").Value = 2.367327375676373
$ws.Range("D8").Value = 18.31797105445999
$ws.Range("E8").Value = 5.349709923656052
$ws.Range("F8").Value = 2.351119530511808
$ws.Range("G8").Value = 4.443672951875615
$ws.Range("H8").Value = 32.82980083617984
$ws.Range("C9").Value = 3.643131942561707
$ws.Range("D9").Value = 19.13816244319345
$ws.Range("E9").Value = 3.706238676419586
$ws.Range("F9").Value = 4.22539330091966
$ws.Range("G9").Value = 5.396128646158712
$ws.Range("H9").Value = 36.10905500925312
$ws.Range("C10").Value = 6.024845510036369
$ws.Range("D10").Value = 18.56090573791954
$ws.Range("E10").Value = 2.822179110313762
$ws.Range("F10").Value = 6.034896197785256
$ws.Range("G10").Value = 5.532599140609939
$ws.Range("H10").Value = 38.97542569666486
